$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new log entry row 48
$ws.Range("A48").Value = 45687
$ws.Range("A48").NumberFormat = $ws.Range("A47").NumberFormat
$ws.Range("B48").Value = "reconfiguring unity stuff"
$ws.Range("C48").Value = 3

# Update selection to reflect where the cursor ended up after entry
$ws.Range("C49").Select()
